$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: columns D1:K1 already exist (restyle not needed, keep s=1) ---
$ws.Range("D1").Value = "CV Train F1"
$ws.Range("E1").Value = "CV Test F1"
$ws.Range("F1").Value = "Validation F1"
$ws.Range("G1").Value = "CV Train Precision"
$ws.Range("H1").Value = "CV Test Precision"
$ws.Range("I1").Value = "Validation Precision"
$ws.Range("J1").Value = "CV Train Recall"
$ws.Range("K1").Value = "CV Test Recall"

# --- Header row: columns L1:O1 are brand new -- copy header formatting first ---
$ws.Range("E1").Copy()
$ws.Range("L1:O1").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("L1").Value = "Validation Recall"
$ws.Range("M1").Value = "Y Val (Validation)"
$ws.Range("N1").Value = "Y Pred (Validation)"
$ws.Range("O1").Value = "Seed"

# --- Row 2: update existing data row (all columns A:O) ---
$ws.Range("A2").Value = "Pipeline(steps=[('scaler', None),`n                ('selector',`n                 <__main__.NamedFeatureSelector object at 0x7fcae05d9850>),`n                ('model',`n                 XGBClassifier(base_score=None, booster=None, callbacks=None,`n                               colsample_bylevel=None, colsample_bynode=None,`n                               colsample_bytree=1.0, early_stopping_rounds=None,`n                               enable_categorical=False, eval_metric=None,`n                               feature_types=None, gamma=0.1, gpu_id=None,`n                               grow_policy=None, importance_type=None,`n                               interaction_constraints=None, learning_rate=0.01,`n                               max_bin=None, max_cat_threshold=None,`n                               max_cat_to_onehot=None, max_delta_step=None,`n                               max_depth=3, max_leaves=None,`n                               min_child_weight=None, missing=nan,`n                               monotone_constraints=None, n_estimators=100,`n                               n_jobs=None, num_parallel_tree=None,`n                               predictor=None, random_state=42, ...))])"
$ws.Range("B2").Value = 0.6326190476190476
$ws.Range("C2").Value = "{'selector': <__main__.NamedFeatureSelector object at 0x7fcae05e2820>, 'scaler': None, 'model__subsample': 1.0, 'model__n_estimators': 100, 'model__max_depth': 3, 'model__learning_rate': 0.01, 'model__gamma': 0.1, 'model__colsample_bytree': 1.0}"
$ws.Range("D2").Value = 0.9762533876637292
$ws.Range("E2").Value = 0.5088329004329003
$ws.Range("F2").Value = 0.7450980392156864
$ws.Range("G2").Value = 0.9659918803802953
$ws.Range("H2").Value = 0.5210507936507937
$ws.Range("I2").Value = 0.6551724137931034
$ws.Range("J2").Value = 0.9878095238095238
$ws.Range("K2").Value = 0.5292
$ws.Range("L2").Value = 0.8636363636363636
$ws.Range("M2").Value = "[1 0 1 1 1 1 0 1 0 1 0 1 0 1 1 0 0 1 1 1 1 0 1 1 0 1 1 1 1 0 0 0 0 1 0 1]"
$ws.Range("N2").Value = "[1 1 1 1 1 1 0 0 1 1 0 1 1 1 0 0 1 1 0 1 1 1 1 1 1 1 1 1 1 1 1 1 0 1 1 1]"
$ws.Range("O2").Value = 42
$ws.Rows.Item(2).AutoFit()

# --- Row 3: update existing data row (all columns A:O) ---
$ws.Range("A3").Value = "Pipeline(steps=[('scaler', RobustScaler()),`n                ('selector',`n                 <__main__.NamedFeatureSelector object at 0x7fcae05e2070>),`n                ('model',`n                 XGBClassifier(base_score=None, booster=None, callbacks=None,`n                               colsample_bylevel=None, colsample_bynode=None,`n                               colsample_bytree=0.5, early_stopping_rounds=None,`n                               enable_categorical=False, eval_metric=None,`n                               feature_types=None, gamma=0, gpu_id=None,`n                               grow_policy=None, importance_type=None,`n                               interaction_constraints=None, learning_rate=0.01,`n                               max_bin=None, max_cat_threshold=None,`n                               max_cat_to_onehot=None, max_delta_step=None,`n                               max_depth=5, max_leaves=None,`n                               min_child_weight=None, missing=nan,`n                               monotone_constraints=None, n_estimators=200,`n                               n_jobs=None, num_parallel_tree=None,`n                               predictor=None, random_state=42, ...))])"
$ws.Range("B3").Value = 0.6852380952380952
$ws.Range("C3").Value = "{'selector': <__main__.NamedFeatureSelector object at 0x7fcae05a8700>, 'scaler': RobustScaler(), 'model__subsample': 0.8, 'model__n_estimators': 200, 'model__max_depth': 5, 'model__learning_rate': 0.01, 'model__gamma': 0, 'model__colsample_bytree': 0.5}"
$ws.Range("D3").Value = 0.964077165417861
$ws.Range("E3").Value = 0.5607665001665002
$ws.Range("F3").Value = 0.6923076923076924
$ws.Range("G3").Value = 0.9506297125349963
$ws.Range("H3").Value = 0.5452904761904762
$ws.Range("I3").Value = 0.6428571428571429
$ws.Range("J3").Value = 0.9793999999999999
$ws.Range("K3").Value = 0.5988
$ws.Range("L3").Value = 0.75
$ws.Range("M3").Value = "[1 1 0 1 0 0 1 0 1 1 1 0 1 1 1 1 1 1 1 1 0 0 1 0 0 1 0 1 1 0 1 1 0 1 1 1]"
$ws.Range("N3").Value = "[0 1 1 1 1 1 1 1 1 1 1 0 1 0 1 1 1 1 1 1 1 1 1 1 1 1 0 0 0 1 0 1 1 1 0 1]"
$ws.Range("O3").Value = 69
$ws.Rows.Item(3).AutoFit()

# --- Row 4: update existing data row (all columns A:O) ---
$ws.Range("A4").Value = "Pipeline(steps=[('scaler', None),`n                ('selector',`n                 <__main__.NamedFeatureSelector object at 0x7fcae05a8430>),`n                ('model',`n                 XGBClassifier(base_score=None, booster=None, callbacks=None,`n                               colsample_bylevel=None, colsample_bynode=None,`n                               colsample_bytree=0.5, early_stopping_rounds=None,`n                               enable_categorical=False, eval_metric=None,`n                               feature_types=None, gamma=0.1, gpu_id=None,`n                               grow_policy=None, importance_type=None,`n                               interaction_constraints=None, learning_rate=0.1,`n                               max_bin=None, max_cat_threshold=None,`n                               max_cat_to_onehot=None, max_delta_step=None,`n                               max_depth=7, max_leaves=None,`n                               min_child_weight=None, missing=nan,`n                               monotone_constraints=None, n_estimators=50,`n                               n_jobs=None, num_parallel_tree=None,`n                               predictor=None, random_state=42, ...))])"
$ws.Range("B4").Value = 0.6233333333333333
$ws.Range("C4").Value = "{'selector': <__main__.NamedFeatureSelector object at 0x7fcae05a8b20>, 'scaler': None, 'model__subsample': 1.0, 'model__n_estimators': 50, 'model__max_depth': 7, 'model__learning_rate': 0.1, 'model__gamma': 0.1, 'model__colsample_bytree': 0.5}"
$ws.Range("D4").Value = 0.9785182521209894
$ws.Range("E4").Value = 0.5576717171717172
$ws.Range("F4").Value = 0.6530612244897959
$ws.Range("G4").Value = 0.9747891661864556
$ws.Range("H4").Value = 0.5521825396825397
$ws.Range("I4").Value = 0.6956521739130435
$ws.Range("J4").Value = 0.982578947368421
$ws.Range("K4").Value = 0.5856
$ws.Range("L4").Value = 0.6153846153846154
$ws.Range("M4").Value = "[0 1 0 0 1 1 1 1 1 1 1 0 1 1 1 1 1 1 1 1 0 1 1 1 0 1 0 1 0 1 0 1 1 1 0 1]"
$ws.Range("N4").Value = "[0 1 1 1 0 1 0 1 1 0 0 1 0 1 1 1 0 1 1 0 1 0 1 1 1 0 0 1 1 1 0 0 1 1 1 1]"
$ws.Range("O4").Value = 23
$ws.Rows.Item(4).AutoFit()

# --- Row 5: update existing data row (all columns A:O) ---
$ws.Range("A5").Value = "Pipeline(steps=[('scaler', RobustScaler()),`n                ('selector',`n                 <__main__.NamedFeatureSelector object at 0x7fcae05a8fd0>),`n                ('model',`n                 XGBClassifier(base_score=None, booster=None, callbacks=None,`n                               colsample_bylevel=None, colsample_bynode=None,`n                               colsample_bytree=0.5, early_stopping_rounds=None,`n                               enable_categorical=False, eval_metric=None,`n                               feature_types=None, gamma=0, gpu_id=None,`n                               grow_policy=None, importance_type=None,`n                               interaction_constraints=None, learning_rate=0.01,`n                               max_bin=None, max_cat_threshold=None,`n                               max_cat_to_onehot=None, max_delta_step=None,`n                               max_depth=3, max_leaves=None,`n                               min_child_weight=None, missing=nan,`n                               monotone_constraints=None, n_estimators=50,`n                               n_jobs=None, num_parallel_tree=None,`n                               predictor=None, random_state=42, ...))])"
$ws.Range("B5").Value = 0.7151190476190477
$ws.Range("C5").Value = "{'selector': <__main__.NamedFeatureSelector object at 0x7fcae05a8fd0>, 'scaler': RobustScaler(), 'model__subsample': 0.5, 'model__n_estimators': 50, 'model__max_depth': 3, 'model__learning_rate': 0.01, 'model__gamma': 0, 'model__colsample_bytree': 0.5}"
$ws.Range("D5").Value = 0.9732271526591363
$ws.Range("E5").Value = 0.5826157842157842
$ws.Range("F5").Value = 0.6382978723404256
$ws.Range("G5").Value = 0.9582585630493827
$ws.Range("H5").Value = 0.5582015873015872
$ws.Range("I5").Value = 0.6
$ws.Range("J5").Value = 0.9908095238095238
$ws.Range("K5").Value = 0.6307999999999998
$ws.Range("L5").Value = 0.6818181818181818
$ws.Range("M5").Value = "[0 1 1 0 0 1 0 0 0 0 1 1 1 0 0 1 1 0 1 1 1 1 1 1 1 1 0 0 1 0 1 1 1 1 1 0]"
$ws.Range("N5").Value = "[0 1 1 1 1 1 0 1 0 1 1 1 0 0 1 0 0 1 1 1 0 1 1 1 1 0 1 1 1 1 0 1 1 0 1 1]"
$ws.Range("O5").Value = 99
$ws.Rows.Item(5).AutoFit()

# --- Row 6: brand new row -- copy formatting (unstyled data cells) from row 5 first ---
$ws.Range("A5:O5").Copy()
$ws.Range("A6:O6").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("A6").Value = "Pipeline(steps=[('scaler', RobustScaler()),`n                ('selector',`n                 <__main__.NamedFeatureSelector object at 0x7fcae05a8eb0>),`n                ('model',`n                 XGBClassifier(base_score=None, booster=None, callbacks=None,`n                               colsample_bylevel=None, colsample_bynode=None,`n                               colsample_bytree=0.8, early_stopping_rounds=None,`n                               enable_categorical=False, eval_metric=None,`n                               feature_types=None, gamma=0.1, gpu_id=None,`n                               grow_policy=None, importance_type=None,`n                               interaction_constraints=None, learning_rate=0.01,`n                               max_bin=None, max_cat_threshold=None,`n                               max_cat_to_onehot=None, max_delta_step=None,`n                               max_depth=5, max_leaves=None,`n                               min_child_weight=None, missing=nan,`n                               monotone_constraints=None, n_estimators=100,`n                               n_jobs=None, num_parallel_tree=None,`n                               predictor=None, random_state=42, ...))])"
$ws.Range("B6").Value = 0.6688888888888889
$ws.Range("C6").Value = "{'selector': <__main__.NamedFeatureSelector object at 0x7fcae05ed6d0>, 'scaler': RobustScaler(), 'model__subsample': 1.0, 'model__n_estimators': 100, 'model__max_depth': 5, 'model__learning_rate': 0.01, 'model__gamma': 0.1, 'model__colsample_bytree': 0.8}"
$ws.Range("D6").Value = 0.9634086897528884
$ws.Range("E6").Value = 0.6123300255300257
$ws.Range("F6").Value = 0.5581395348837209
$ws.Range("G6").Value = 0.9430689646559474
$ws.Range("H6").Value = 0.5559349206349207
$ws.Range("I6").Value = 0.5217391304347826
$ws.Range("J6").Value = 0.9872727272727272
$ws.Range("K6").Value = 0.6988000000000001
$ws.Range("L6").Value = 0.6
$ws.Range("M6").Value = "[1 0 1 1 0 0 0 0 1 0 1 1 0 1 1 0 1 0 0 0 0 0 1 1 1 0 1 0 1 1 1 1 1 1 1 0]"
$ws.Range("N6").Value = "[1 1 1 1 1 1 1 0 0 0 0 0 1 1 1 1 0 1 1 0 1 1 1 1 1 1 1 0 1 1 0 1 0 0 0 0]"
$ws.Range("O6").Value = 89
$ws.Rows.Item(6).AutoFit()
